# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The data table (rows 16-52, columns E:G) lists 37 monthly "Periodo Mora"
# entries. They used to run in descending order (2003 down to 1703); now
# they run in ascending chronological order (1703 up to 2003), and the
# "Valor Mora" (F) / "Salario Basico" (G) figures for each period are
# refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 52

# Starting period (YYMM, two-digit year + two-digit month) and how many
# of the 37 rows keep the "old" F value (29509) before switching to the
# "new" F value (31249). Periods 1703-1808 -> 29509, 1809-2003 -> 31249.
$year  = 17
$month = 3

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $period = ("{0:D2}{1:D2}" -f $year, $month)

    $ws.Range("E$row").Value = $period

    if ([int]$period -le 1808) {
        $ws.Range("F$row").Value = 29509
    } else {
        $ws.Range("F$row").Value = 31249
    }

    $ws.Range("G$row").Value = 781242

    $month = $month + 1
    if ($month -gt 12) {
        $month = 1
        $year = $year + 1
    }
}
